$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 4). Set cells in an order that reproduces the
# exact shared-string insertion order seen in the target workbook:
# MABROUK, WAFAE, YOUSSEFI, ABDELMALEK (new strings), reusing existing
# shared strings for CONSULTANT / BCP / HR.
$ws.Range("H4").Value = "MABROUK"
$ws.Range("I4").Value = "WAFAE"
$ws.Range("B4").Value = "YOUSSEFI"
$ws.Range("C4").Value = "ABDELMALEK"

$ws.Range("A4").Value = 7573
$ws.Range("D4").Value = "CONSULTANT"
$ws.Range("E4").Value = "BCP"
$ws.Range("F4").Value = "HR"
$ws.Range("G4").Value = 7373

# J4 is a date (BIRTHDAY MANAGER N+1). Copy the number format from J2 so the
# new cell reuses the existing date style instead of creating a new one,
# then set its value.
$ws.Range("J2").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = 35802

# Update the sheet view: clear the frozen/scrolled topLeftCell (scroll back
# to the default top-left position) and move the selection to E8.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E8").Select() | Out-Null
